$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "ORG_GOV_IDENOLD"
$ws.Range("D1").Value = "ORG_GOV_IDENNEW"
$ws.Range("E1").Value = "ORG_GOV_STATUS"

$ws.Range("E2").Select()
